# Add a third worksheet "Find" to the workbook by copying the existing
# "Create" sheet (same headers/data: Uname, pwd, companyName, firstName,
# lastName), placing it after "Create". Also update the selection state
# on both sheets to match the post-edit UI state (Find becomes the
# active/tab-selected sheet).

$wb = $excel.ActiveWorkbook

$createSheet = $wb.Worksheets.Item("Create")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate "Create" (carries over data, styles, column widths) and place
# the copy right after the last existing sheet, then rename it to "Find".
$createSheet.Copy($null, $lastSheet)
$findSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$findSheet.Name = "Find"

# Update the selection left behind on "Create" (no longer the active tab).
$createSheet.Activate()
$createSheet.Range("A1:E3").Select()

# "Find" becomes the active / tab-selected sheet, with its own selection.
$findSheet.Activate()
$findSheet.Range("D6").Select()
